$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.198.36'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '2.279.35'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +1.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.34'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.31'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.621'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.26%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.594'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.54'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0896'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.20'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.35%  '
$ws.Range("E13").Value = '  +1.63%  '
$ws.Range("E14").Value = '  +1.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.97'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.76%  '
$ws.Range("D16").Value = '2.624.22'
$ws.Range("E16").Value = '  -0.73%  '
$ws.Range("D17").Value = '2.323.06'
$ws.Range("E17").Value = '  +1.12%  '
$ws.Range("D18").Value = '42.325.45'
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("E19").Value = '  -2.70%  '
$ws.Range("E20").Value = '  -1.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.45'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +8.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.80'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.81%  '
$ws.Range("E23").Value = '  -2.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '263.38'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -4.37%  '
$ws.Range("E25").Value = '  -4.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.01'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.48%  '
$ws.Range("E27").Value = '  -2.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.34'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.83'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +12.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.26'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.74'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -4.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '164.47'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.91%  '
$ws.Range("E33").Value = '  -2.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.130'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.37%  '
$ws.Range("E35").Value = '  -1.66%  '
$ws.Range("E36").Value = '  -5.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.44'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0346'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.68'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.63'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.43%  '
$ws.Range("E41").Value = '  +4.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.90'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.74%  '
$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '68.64'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.63%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.225'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '11.86'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.81%  '
$ws.Range("D47").Value = '1.698.76'
$ws.Range("E47").Value = '  +6.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '78.99'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '109.64'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.62'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.13'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.66%  '
